# Auto-generated Excel COM-interop script
# Applies the numeric corrections to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets
# as captured by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 882.2239
$ws.Range("J17").Value = 882.2239
$ws.Range("L17").Value = 2646.6717
$ws.Range("N17").Value = -2982.6717

$ws.Range("H69").Value = 3220
$ws.Range("J69").Value = 4900
$ws.Range("L69").Value = 14700
$ws.Range("N69").Value = -16448

$ws.Range("H72").Value = 3220
$ws.Range("J72").Value = 4900
$ws.Range("L72").Value = 44100
$ws.Range("N72").Value = -52836

$ws.Range("H76").Value = 3351321.2
$ws.Range("I76").Value = 4688330.5
$ws.Range("J76").Value = 8797.5
$ws.Range("K76").Value = 4688330.5
$ws.Range("L76").Value = 8797.5
$ws.Range("M76").Value = -4688015.5
$ws.Range("N76").Value = -9427.5

$ws.Range("H79").Value = 3351321.2
$ws.Range("I79").Value = 4688330.5
$ws.Range("J79").Value = 8797.5
$ws.Range("K79").Value = 4688330.5
$ws.Range("L79").Value = 8797.5
$ws.Range("M79").Value = -4687238.5
$ws.Range("N79").Value = -10981.5

$ws.Range("H80").Value = 1162.4736
$ws.Range("I80").Value = 1760.1
$ws.Range("K80").Value = 5280.299999999999
$ws.Range("M80").Value = -4282.299999999999

$ws.Range("H83").Value = 1162.4736
$ws.Range("I83").Value = 1760.1
$ws.Range("K83").Value = 15840.9
$ws.Range("M83").Value = -10848.9

$ws.Range("H112").Value = 4722.154
$ws.Range("J112").Value = 4722.154
$ws.Range("L112").Value = 14166.462
$ws.Range("N112").Value = -16382.462

$ws.Range("H138").Value = 2748.9272
$ws.Range("I138").Value = 2358.0334
$ws.Range("K138").Value = 7074.100199999999
$ws.Range("M138").Value = -1934.100199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3686.349
$ws.Range("I32").Value = 3017.0613
$ws.Range("K32").Value = 3017.0613
$ws.Range("M32").Value = -2730.0613

$ws.Range("H34").Value = 10000
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10542

$ws.Range("H132").Value = 2136.0386
$ws.Range("I132").Value = 1710.6428
$ws.Range("J132").Value = 2632.3333
$ws.Range("K132").Value = 5131.928400000001
$ws.Range("L132").Value = 7896.999899999999
$ws.Range("M132").Value = -2601.928400000001
$ws.Range("N132").Value = -12956.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 714.3333
$ws.Range("I64").Value = 717.2
$ws.Range("J64").Value = 700
$ws.Range("K64").Value = 717.2
$ws.Range("L64").Value = 700
$ws.Range("M64").Value = -492.2
$ws.Range("N64").Value = -1150

$ws.Range("H67").Value = 714.3333
$ws.Range("I67").Value = 717.2
$ws.Range("J67").Value = 700
$ws.Range("K67").Value = 717.2
$ws.Range("L67").Value = 700
$ws.Range("M67").Value = 62.79999999999995
$ws.Range("N67").Value = -2260

$ws.Range("H94").Value = 338.97058
$ws.Range("I94").Value = 300.51614
$ws.Range("K94").Value = 300.51614
$ws.Range("M94").Value = 150.48386

$ws.Range("H134").Value = 8644.15
$ws.Range("I134").Value = 10305.8125
$ws.Range("J134").Value = 1997.5
$ws.Range("K134").Value = 30917.4375
$ws.Range("L134").Value = 5992.5
$ws.Range("M134").Value = -28382.4375
$ws.Range("N134").Value = -11062.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1747.3077
$ws.Range("I31").Value = 1034
$ws.Range("J31").Value = 2270.4
$ws.Range("K31").Value = 1034
$ws.Range("L31").Value = 2270.4
$ws.Range("M31").Value = -739
$ws.Range("N31").Value = -2860.4

$ws.Range("H34").Value = 1747.3077
$ws.Range("I34").Value = 1034
$ws.Range("J34").Value = 2270.4
$ws.Range("K34").Value = 1034
$ws.Range("L34").Value = 2270.4
$ws.Range("M34").Value = -832
$ws.Range("N34").Value = -2674.4

$ws.Range("H106").Value = 34825
$ws.Range("J106").Value = 34650
$ws.Range("L106").Value = 34650
$ws.Range("N106").Value = -37174

$ws.Range("H132").Value = 2201.8333
$ws.Range("I132").Value = 1025.4706
$ws.Range("J132").Value = 5058.7144
$ws.Range("K132").Value = 3076.4118
$ws.Range("L132").Value = 15176.1432
$ws.Range("M132").Value = -546.4118000000003
$ws.Range("N132").Value = -20236.1432

$ws.Range("H134").Value = 1379.3
$ws.Range("I134").Value = 1388.7368
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 4166.2104
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -1631.2104
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 817.5
$ws.Range("I3").Value = 817.5
$ws.Range("K3").Value = 2452.5
$ws.Range("M3").Value = -2340.5

$ws.Range("H5").Value = 812.7778
$ws.Range("J5").Value = 835.8333
$ws.Range("L5").Value = 2507.4999
$ws.Range("N5").Value = -2731.4999

$ws.Range("H33").Value = 152.33333
$ws.Range("I33").Value = 172.33333
$ws.Range("J33").Value = 132.33333
$ws.Range("K33").Value = 1033.99998
$ws.Range("L33").Value = 793.9999799999999
$ws.Range("M33").Value = -750.9999800000001
$ws.Range("N33").Value = -1359.99998

$ws.Range("H121").Value = 789.0909
$ws.Range("I121").Value = 815
$ws.Range("J121").Value = 783.3333
$ws.Range("K121").Value = 2445
$ws.Range("L121").Value = 2349.9999
$ws.Range("M121").Value = -1135
$ws.Range("N121").Value = -4969.9999

$ws.Range("H122").Value = 932.9167
$ws.Range("I122").Value = 523
$ws.Range("J122").Value = 1225.7142
$ws.Range("K122").Value = 4707
$ws.Range("L122").Value = 11031.4278
$ws.Range("M122").Value = -2257
$ws.Range("N122").Value = -15931.4278

$ws.Range("H131").Value = 783.47
$ws.Range("J131").Value = 794.23956
$ws.Range("L131").Value = 2382.71868
$ws.Range("N131").Value = -12462.71868

$ws.Range("H135").Value = 812.7778
$ws.Range("J135").Value = 835.8333
$ws.Range("L135").Value = 7522.4997
$ws.Range("N135").Value = -12592.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2490.16
$ws.Range("I102").Value = 2496.842
$ws.Range("J102").Value = 2469
$ws.Range("K102").Value = 2496.842
$ws.Range("L102").Value = 2469
$ws.Range("M102").Value = -874.8420000000001
$ws.Range("N102").Value = -5713

$ws.Range("H132").Value = 2749973.5
$ws.Range("J132").Value = 3596.8572
$ws.Range("L132").Value = 10790.5716
$ws.Range("N132").Value = -15850.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2590.5
$ws.Range("I16").Value = 2499.1428
$ws.Range("J16").Value = 2803.6667
$ws.Range("K16").Value = 2499.1428
$ws.Range("L16").Value = 2803.6667
$ws.Range("M16").Value = -2329.1428
$ws.Range("N16").Value = -3143.6667

$ws.Range("H40").Value = 4494
$ws.Range("I40").Value = 1628.7273
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 1628.7273
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -1492.7273
$ws.Range("N40").Value = -15272

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H107").Value = 590.5333000000001
$ws.Range("J107").Value = 731.2857
$ws.Range("L107").Value = 2193.8571
$ws.Range("N107").Value = -6033.8571

$ws.Range("H132").Value = 1105.625
$ws.Range("I132").Value = 936.34784
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 2809.04352
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -279.0435200000002
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 27780980
$ws.Range("I136").Value = 39685950
$ws.Range("K136").Value = 119057850
$ws.Range("M136").Value = -119055300
